# Apply scheduled-runner price/profit updates to the Garuda_Profits workbook
# Each block below targets one (sheet, row) pair; values come from the latest price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 10
$ws.Range("H10").Value = 1934.6666
$ws.Range("I10").Value = 2004
$ws.Range("J10").Value = 1900
$ws.Range("K10").Value = 2004
$ws.Range("L10").Value = 1900
$ws.Range("M10").Value = -1711
$ws.Range("N10").Value = -2486

# Row 40
$ws.Range("H40").Value = 2180
$ws.Range("I40").Value = 2859.875
$ws.Range("J40").Value = 1761.6154
$ws.Range("K40").Value = 2859.875
$ws.Range("L40").Value = 1761.6154
$ws.Range("M40").Value = -2684.875
$ws.Range("N40").Value = -2111.6154

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 806.07275
$ws.Range("I74").Value = 815.2075
$ws.Range("J74").Value = 564
$ws.Range("K74").Value = 815.2075
$ws.Range("L74").Value = 564
$ws.Range("M74").Value = 58.79250000000002
$ws.Range("N74").Value = -2312

# Row 77
$ws.Range("H77").Value = 806.07275
$ws.Range("I77").Value = 815.2075
$ws.Range("J77").Value = 564
$ws.Range("K77").Value = 4076.0375
$ws.Range("L77").Value = 2820
$ws.Range("M77").Value = 291.9625000000001
$ws.Range("N77").Value = -11556

# Row 88
$ws.Range("H88").Value = 10922.6
$ws.Range("I88").Value = 11503
$ws.Range("J88").Value = 10535.667
$ws.Range("K88").Value = 11503
$ws.Range("L88").Value = 10535.667
$ws.Range("M88").Value = -11097
$ws.Range("N88").Value = -11347.667

# Row 91
$ws.Range("H91").Value = 10922.6
$ws.Range("I91").Value = 11503
$ws.Range("J91").Value = 10535.667
$ws.Range("K91").Value = 11503
$ws.Range("L91").Value = 10535.667
$ws.Range("M91").Value = -10099
$ws.Range("N91").Value = -13343.667

$ws = $wb.Worksheets.Item("BSM")
# Row 12
$ws.Range("H12").Value = 399.66666
$ws.Range("I12").Value = 399.66666
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 399.66666
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -231.66666
$ws.Range("N12").ClearContents()

# Row 14
$ws.Range("H14").Value = 12000
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

# Row 64
$ws.Range("H64").Value = 461.75
$ws.Range("I64").Value = 430.8
$ws.Range("K64").Value = 430.8
$ws.Range("M64").Value = -205.8

# Row 67
$ws.Range("H67").Value = 461.75
$ws.Range("I67").Value = 430.8
$ws.Range("K67").Value = 430.8
$ws.Range("M67").Value = 349.2

# Row 134
$ws.Range("H134").Value = 3655.7234
$ws.Range("I134").Value = 3671.1135
$ws.Range("J134").Value = 3430
$ws.Range("K134").Value = 11013.3405
$ws.Range("L134").Value = 10290
$ws.Range("M134").Value = -8478.3405
$ws.Range("N134").Value = -15360

$ws = $wb.Worksheets.Item("CRP")
# Row 52
$ws.Range("H52").Value = 93850
$ws.Range("J52").Value = 93850
$ws.Range("L52").Value = 93850
$ws.Range("N52").Value = -94438

# Row 58
$ws.Range("H58").Value = 1383.7333
$ws.Range("I58").Value = 1383.7333
$ws.Range("K58").Value = 1383.7333
$ws.Range("M58").Value = -1180.7333

# Row 99
$ws.Range("H99").Value = 2099.087
$ws.Range("I99").Value = 2016.3636
$ws.Range("J99").Value = 2174.9167
$ws.Range("K99").Value = 2016.3636
$ws.Range("L99").Value = 2174.9167
$ws.Range("M99").Value = -518.3635999999999
$ws.Range("N99").Value = -5170.9167

# Row 126
$ws.Range("H126").Value = 2099.087
$ws.Range("I126").Value = 2016.3636
$ws.Range("J126").Value = 2174.9167
$ws.Range("K126").Value = 6049.0908
$ws.Range("L126").Value = 6524.750100000001
$ws.Range("M126").Value = -3579.0908
$ws.Range("N126").Value = -11464.7501

# Row 136
$ws.Range("H136").Value = 1383.7333
$ws.Range("I136").Value = 1383.7333
$ws.Range("K136").Value = 4151.199900000001
$ws.Range("M136").Value = -1601.199900000001

$ws = $wb.Worksheets.Item("CUL")
# Row 10
$ws.Range("H10").Value = 110.5
$ws.Range("I10").Value = 110.5
$ws.Range("K10").Value = 331.5
$ws.Range("M10").Value = -192.5

$ws = $wb.Worksheets.Item("GSM")
# Row 17
$ws.Range("H17").Value = 60
$ws.Range("I17").Value = 60
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 60
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 108
$ws.Range("N17").ClearContents()

# Row 80
$ws.Range("H80").Value = 4306.524
$ws.Range("J80").Value = 4001.4666
$ws.Range("L80").Value = 4001.4666
$ws.Range("N80").Value = -5997.4666

# Row 83
$ws.Range("H83").Value = 4306.524
$ws.Range("J83").Value = 4001.4666
$ws.Range("L83").Value = 20007.333
$ws.Range("N83").Value = -29991.333

# Row 139
$ws.Range("H139").Value = 37812.6
$ws.Range("J139").Value = 37812.6
$ws.Range("L139").Value = 37812.6
$ws.Range("N139").Value = -48092.6

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 38463984
$ws.Range("I7").Value = 2333.3333
$ws.Range("K7").Value = 2333.3333
$ws.Range("M7").Value = -2221.3333

# Row 48
$ws.Range("H48").Value = 29997.5
$ws.Range("I48").Value = 3000
$ws.Range("J48").Value = 38996.668
$ws.Range("K48").Value = 3000
$ws.Range("L48").Value = 38996.668
$ws.Range("M48").Value = -2339
$ws.Range("N48").Value = -40318.668

# Row 122
$ws.Range("H122").Value = 4807.357
$ws.Range("I122").Value = 5391.1816
$ws.Range("J122").Value = 2666.6667
$ws.Range("K122").Value = 16173.5448
$ws.Range("L122").Value = 8000.000100000001
$ws.Range("M122").Value = -13723.5448
$ws.Range("N122").Value = -12900.0001

# Row 126
$ws.Range("H126").Value = 38463984
$ws.Range("I126").Value = 2333.3333
$ws.Range("K126").Value = 6999.999899999999
$ws.Range("M126").Value = -4529.999899999999

# Row 132
$ws.Range("H132").Value = 7718.343
$ws.Range("I132").Value = 10529.772
$ws.Range("J132").Value = 2960.5386
$ws.Range("K132").Value = 31589.316
$ws.Range("L132").Value = 8881.6158
$ws.Range("M132").Value = -29059.316
$ws.Range("N132").Value = -13941.6158

# Row 136
$ws.Range("H136").Value = 13069
$ws.Range("I136").Value = 15461.25
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 46383.75
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -43833.75
$ws.Range("N136").Value = -15600

$ws = $wb.Worksheets.Item("WVR")
# Row 19
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

# Row 107
$ws.Range("H107").Value = 15804976
$ws.Range("I107").Value = 7353374
$ws.Range("J107").Value = 27778078
$ws.Range("K107").Value = 22060122
$ws.Range("L107").Value = 83334234
$ws.Range("M107").Value = -22058202
$ws.Range("N107").Value = -83338074

# Row 122
$ws.Range("H122").Value = 1713.9474
$ws.Range("I122").Value = 1624.909
$ws.Range("J122").Value = 2015.3077
$ws.Range("K122").Value = 4874.727000000001
$ws.Range("L122").Value = 6045.9231
$ws.Range("M122").Value = -2424.727000000001
$ws.Range("N122").Value = -10945.9231

# Row 136
$ws.Range("H136").Value = 1329.1187
$ws.Range("I136").Value = 1193.2858
$ws.Range("J136").Value = 1664.7059
$ws.Range("K136").Value = 3579.8574
$ws.Range("L136").Value = 4994.1177
$ws.Range("M136").Value = -1029.8574
$ws.Range("N136").Value = -10094.1177
